# Natmi following Dr Hou advice
# Recompute the Reln -> Itga3 ligand/receptor stats now that the clustering adds
# the FAPs -> FAPs combination. The table now covers every Sending x Target pair
# across the 3 clusters (ECs, FAPs, sCs) = 9 rows instead of 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A="ECs"; D="ECs"; E=3; F=1; G=0.03565433333333334; H=0.106963; I=0.002412342638581826; J=0.002412342638581825; K=3; L=1; M=8.269168666666666; N=24.807506; O=0.671680253471746; P=0.671680253471746; Q=0.2948316960308889; R=2.653485264278; S=0.001620322914943341; T=0.001620322914943341 }
    @{ Row=3; A="ECs"; D="FAPs"; E=3; F=1; G=0.03565433333333334; H=0.106963; I=0.002412342638581826; J=0.002412342638581825; K=2; L=0.6666666666666666; M=0.121294; N=0.363882; O=0.009852354928133683; P=0.009852354928133683; Q=0.004324656707333334; R=0.038921910366; S=0.00002376725588357867; T=0.00002376725588357866 }
    @{ Row=4; A="ECs"; D="sCs"; E=3; F=1; G=0.03565433333333334; H=0.106963; I=0.002412342638581826; J=0.002412342638581825; K=3; L=1; M=3.920705666666667; N=11.762117; O=0.3184673916001203; P=0.3184673916001203; Q=0.1397901467412222; R=1.258111320671; S=0.0007682524677549057; T=0.0007682524677549056 }
    @{ Row=5; A="FAPs"; D="ECs"; E=3; F=1; G=8.066615333333333; H=24.199846; I=0.5457805068380079; J=0.5457805068380079; K=3; L=1; M=8.269168666666666; N=24.807506; O=0.671680253471746; P=0.671680253471746; Q=66.70420276045287; R=600.337824844076; S=0.3665899891728912; T=0.3665899891728912 }
    @{ Row=6; A="FAPs"; D="FAPs"; E=3; F=1; G=8.066615333333333; H=24.199846; I=0.5457805068380079; J=0.5457805068380079; K=2; L=0.6666666666666666; M=0.121294; N=0.363882; O=0.009852354928133683; P=0.009852354928133683; Q=0.9784320402413332; R=8.805888362172; S=0.005377223266224747; T=0.005377223266224747 }
    @{ Row=7; A="FAPs"; D="sCs"; E=3; F=1; G=8.066615333333333; H=24.199846; I=0.5457805068380079; J=0.5457805068380079; K=3; L=1; M=3.920705666666667; N=11.762117; O=0.3184673916001203; P=0.3184673916001203; Q=31.62682444822022; R=284.641420033982; S=0.173813294398892; T=0.173813294398892 }
    @{ Row=8; A="sCs"; D="ECs"; E=3; F=1; G=6.677692666666666; H=20.033078; I=0.4518071505234102; J=0.4518071505234102; K=3; L=1; M=8.269168666666666; N=24.807506; O=0.671680253471746; P=0.671680253471746; Q=55.21896696482977; R=496.9707026834679; S=0.3034699413839114; T=0.3034699413839114 }
    @{ Row=9; A="sCs"; D="FAPs"; E=3; F=1; G=6.677692666666666; H=20.033078; I=0.4518071505234102; J=0.4518071505234102; K=2; L=0.6666666666666666; M=0.121294; N=0.363882; O=0.009852354928133683; P=0.009852354928133683; Q=0.8099640543106666; R=7.289676488795999; S=0.004451364406025358; T=0.004451364406025358 }
    @{ Row=10; A="sCs"; D="sCs"; E=3; F=1; G=6.677692666666666; H=20.033078; I=0.4518071505234102; J=0.4518071505234102; K=3; L=1; M=3.920705666666667; N=11.762117; O=0.3184673916001203; P=0.3184673916001203; Q=26.18126747845844; R=235.631407306126; S=0.1438858447334734; T=0.1438858447334734 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = "Reln"
    $ws.Range("C$row").Value = "Itga3"
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
}
